# slides: minor slide edits
# Add speaker notes to slide 16 ("Aliases") explaining how `python filename.py`
# is used to run a Python script from a text file.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# Touching NotesPage ensures the notes part/placeholder machinery is ready;
# the body (speaker notes) placeholder is added/fetched explicitly since this
# host only materializes the notes body placeholder on demand.
$notesPage = $s.NotesPage
$notesBody = $notesPage.Shapes.AddPlaceholder(2)

$notesBody.TextFrame.TextRange.Text = "For those who may not have seen it before, “python filename.py” is how you would run a text file containing python code. "
